$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 5991.4287
$ws.Range("I13").Value = 4666.3335
$ws.Range("J13").Value = 6985.25
$ws.Range("K13").Value = 4666.3335
$ws.Range("L13").Value = 6985.25
$ws.Range("M13").Value = -4497.3335
$ws.Range("N13").Value = -7323.25
$ws.Range("H33").Value = 385253.2
$ws.Range("I33").Value = 526567.9
$ws.Range("K33").Value = 526567.9
$ws.Range("M33").Value = -526338.9
$ws.Range("H51").Value = 3333.1667
$ws.Range("J51").Value = 3399.8
$ws.Range("L51").Value = 3399.8
$ws.Range("N51").Value = -4367.8
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H106").Value = 9628.182000000001
$ws.Range("I106").Value = 8948.5
$ws.Range("K106").Value = 8948.5
$ws.Range("M106").Value = -8317.5
$ws.Range("H138").Value = 5385.915
$ws.Range("J138").Value = 5431.698
$ws.Range("L138").Value = 16295.094
$ws.Range("N138").Value = -26575.094

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 16668267
$ws.Range("J13").Value = 4750
$ws.Range("L13").Value = 4750
$ws.Range("N13").Value = -5038
$ws.Range("H32").Value = 2465.831
$ws.Range("I32").Value = 2009.2192
$ws.Range("K32").Value = 2009.2192
$ws.Range("M32").Value = -1722.2192
$ws.Range("H46").Value = 9810
$ws.Range("J46").Value = 8512.5
$ws.Range("L46").Value = 8512.5
$ws.Range("N46").Value = -9150.5
$ws.Range("H74").Value = 1552.1034
$ws.Range("I74").Value = 1430.0358
$ws.Range("K74").Value = 1430.0358
$ws.Range("M74").Value = -556.0358000000001
$ws.Range("H77").Value = 1552.1034
$ws.Range("I77").Value = 1430.0358
$ws.Range("K77").Value = 7150.179
$ws.Range("M77").Value = -2782.179
$ws.Range("H86").Value = 39999
$ws.Range("J86").Value = 39999
$ws.Range("L86").Value = 39999
$ws.Range("N86").Value = -42371
$ws.Range("H89").Value = 39999
$ws.Range("J89").Value = 39999
$ws.Range("L89").Value = 119997
$ws.Range("N89").Value = -131853
$ws.Range("H122").Value = 4097.4614
$ws.Range("I122").Value = 4026.7
$ws.Range("K122").Value = 12080.1
$ws.Range("M122").Value = -9630.099999999999
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 26332.666
$ws.Range("J57").Value = 26332.666
$ws.Range("L57").Value = 26332.666
$ws.Range("N57").Value = -27772.666
$ws.Range("H96").Value = 23216.75
$ws.Range("I96").Value = 8456.5
$ws.Range("J96").Value = 67497.5
$ws.Range("K96").Value = 8456.5
$ws.Range("L96").Value = 67497.5
$ws.Range("M96").Value = -5710.5
$ws.Range("N96").Value = -72989.5
$ws.Range("H105").Value = 5043.5
$ws.Range("I105").Value = 4384.933
$ws.Range("K105").Value = 4384.933
$ws.Range("M105").Value = -2637.933
$ws.Range("H107").Value = 5614.5454
$ws.Range("I107").Value = 3708.8462
$ws.Range("K107").Value = 3708.8462
$ws.Range("M107").Value = -1788.8462
$ws.Range("H120").Value = 70000
$ws.Range("J120").Value = 70000
$ws.Range("L120").Value = 70000
$ws.Range("N120").Value = -79676
$ws.Range("H132").Value = 296000
$ws.Range("J132").Value = 296000
$ws.Range("L132").Value = 296000
$ws.Range("N132").Value = -306120
$ws.Range("H134").Value = 2762.9614
$ws.Range("I134").Value = 2841.4211
$ws.Range("K134").Value = 8524.263300000001
$ws.Range("M134").Value = -5989.263300000001
$ws.Range("H136").Value = 26332.666
$ws.Range("J136").Value = 26332.666
$ws.Range("L136").Value = 26332.666
$ws.Range("N136").Value = -36532.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5114.737
$ws.Range("I31").Value = 2306.7273
$ws.Range("K31").Value = 2306.7273
$ws.Range("M31").Value = -2011.7273
$ws.Range("H34").Value = 5114.737
$ws.Range("I34").Value = 2306.7273
$ws.Range("K34").Value = 2306.7273
$ws.Range("M34").Value = -2104.7273
$ws.Range("H52").Value = 49850
$ws.Range("J52").Value = 49850
$ws.Range("L52").Value = 49850
$ws.Range("N52").Value = -50438
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H134").Value = 2279.5356
$ws.Range("I134").Value = 2175.4546
$ws.Range("K134").Value = 6526.3638
$ws.Range("M134").Value = -3991.3638
$ws.Range("H138").Value = 382894.75
$ws.Range("J138").Value = 382894.75
$ws.Range("L138").Value = 382894.75
$ws.Range("N138").Value = -393174.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 6554.5
$ws.Range("I7").Value = 386
$ws.Range("K7").Value = 1158
$ws.Range("M7").Value = -1046
$ws.Range("H34").Value = 3048.75
$ws.Range("I34").Value = 481.8889
$ws.Range("K34").Value = 1445.6667
$ws.Range("M34").Value = -1361.6667
$ws.Range("H112").Value = 3767
$ws.Range("I112").Value = 2834.75
$ws.Range("J112").Value = 5010
$ws.Range("K112").Value = 8504.25
$ws.Range("L112").Value = 15030
$ws.Range("M112").Value = -7396.25
$ws.Range("N112").Value = -17246
$ws.Range("H113").Value = 1972.4
$ws.Range("J113").Value = 2578.5
$ws.Range("L113").Value = 7735.5
$ws.Range("N113").Value = -12075.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6707.1113
$ws.Range("I80").Value = 5455.3335
$ws.Range("J80").Value = 7333
$ws.Range("K80").Value = 5455.3335
$ws.Range("L80").Value = 7333
$ws.Range("M80").Value = -4457.3335
$ws.Range("N80").Value = -9329
$ws.Range("H83").Value = 6707.1113
$ws.Range("I83").Value = 5455.3335
$ws.Range("J83").Value = 7333
$ws.Range("K83").Value = 27276.6675
$ws.Range("L83").Value = 36665
$ws.Range("M83").Value = -22284.6675
$ws.Range("N83").Value = -46649
$ws.Range("H102").Value = 14200
$ws.Range("I102").Value = 3666.6667
$ws.Range("K102").Value = 3666.6667
$ws.Range("M102").Value = -2044.6667
$ws.Range("H122").Value = 14855
$ws.Range("I122").Value = 15548.571
$ws.Range("K122").Value = 46645.713
$ws.Range("M122").Value = -44195.713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2598.6924
$ws.Range("J46").Value = 3161.111
$ws.Range("L46").Value = 3161.111
$ws.Range("N46").Value = -3537.111
$ws.Range("H55").Value = 1772.909
$ws.Range("J55").Value = 3368.7
$ws.Range("L55").Value = 3368.7
$ws.Range("N55").Value = -3714.7
$ws.Range("I74").Value = 44997.5
$ws.Range("K74").Value = 44997.5
$ws.Range("M74").Value = -43999.5
$ws.Range("I77").Value = 44997.5
$ws.Range("K77").Value = 134992.5
$ws.Range("M77").Value = -130000.5
$ws.Range("H82").Value = 2970.6191
$ws.Range("J82").Value = 3322.353
$ws.Range("L82").Value = 3322.353
$ws.Range("N82").Value = -4044.353
$ws.Range("H85").Value = 2970.6191
$ws.Range("J85").Value = 3322.353
$ws.Range("L85").Value = 3322.353
$ws.Range("N85").Value = -5818.353

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2462.2856
$ws.Range("I81").Value = 1206
$ws.Range("J81").Value = 10000
$ws.Range("K81").Value = 2412
$ws.Range("L81").Value = 20000
$ws.Range("M81").Value = -1351
$ws.Range("N81").Value = -22122
$ws.Range("H84").Value = 2462.2856
$ws.Range("I84").Value = 1206
$ws.Range("J84").Value = 10000
$ws.Range("K84").Value = 12060
$ws.Range("L84").Value = 100000
$ws.Range("M84").Value = -6756
$ws.Range("N84").Value = -110608
$ws.Range("H107").Value = 5199.875
$ws.Range("I107").Value = 4809.8
$ws.Range("K107").Value = 14429.4
$ws.Range("M107").Value = -12509.4
$ws.Range("H132").Value = 2866.3076
$ws.Range("I132").Value = 1442.1111
$ws.Range("J132").Value = 4087.0476
$ws.Range("K132").Value = 4326.3333
$ws.Range("L132").Value = 12261.1428
$ws.Range("M132").Value = -1796.3333
$ws.Range("N132").Value = -17321.1428
